$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add materials for session 11 (week 11 / row 12) and prep for week 12 (row 13)
$ws.Range("F12").Value = "exercises/e11.html"
$ws.Range("D13").Value = "prep/p12.html"
$ws.Range("E12").Value = "slides/slides.html#/sitzung-11-machine-agency-wie-algorithmen-das-unterhaltungserleben-beeiflussen-können"

# Update the active selection to match the edited workbook state
$ws.Range("E12").Select()
